# river update May 2024
# Applies the May-2024 refresh of trend-result values for
# "Whitebait Creek at Edinburgh Terrace" (Rivers trend results sheet).
#
# Summary of the change:
#  - Rows 2-10 (Visual Clarity .. SIN) keep their parameter/site columns but get
#    refreshed statistic values (F,G,H,I,J,K,L,M,N) and refreshed confidence
#    wording in column P (and a couple of D/"seasonal trend" flips).
#  - The old row 11 ("Suspended Sediment Concentration") drops out of this
#    year's results entirely; "Total Nitrogen" (previously row 12) now occupies
#    row 11 with refreshed values, and "Total Phosphorus" (previously row 13)
#    now occupies row 12 with refreshed values.
#  - The sheet therefore shrinks from 13 data rows to 12 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("F2").Value = 0.993518344364968
    $ws.Range("H2").Value = 0.833333333333333
    $ws.Range("J2").Value = 0.365
    $ws.Range("K2").Value = 0.0406661295737444
    $ws.Range("L2").Value = 0.0152881391792108
    $ws.Range("M2").Value = 0.0653063317257669
    $ws.Range("N2").Value = 11.1414053626697
    $ws.Range("P2").Value = 'Virtually certain improving'

    # Row 3
    $ws.Range("F3").Value = 0.643788875060667
    $ws.Range("J3").Value = 8.28
    $ws.Range("K3").Value = 0.153922804922506
    $ws.Range("L3").Value = -0.347680587409719
    $ws.Range("M3").Value = 0.359472331984342
    $ws.Range("N3").Value = 1.8589710739433
    $ws.Range("P3").Value = 'As likely as not increasing'

    # Row 4
    $ws.Range("F4").Value = 0.0449824875182887
    $ws.Range("H4").Value = 0.913793103448276
    $ws.Range("J4").Value = 0.0785
    $ws.Range("K4").Value = 0.005103395759259
    $ws.Range("L4").Value = 0.0001928836505708
    $ws.Range("M4").Value = 0.0105515059183995
    $ws.Range("N4").Value = 6.50114109459755
    $ws.Range("P4").Value = 'Extremely unlikely improving'

    # Row 5
    $ws.Range("F5").Value = 0.872963383828828
    $ws.Range("G5").Value = 0.0172413793103448
    $ws.Range("H5").Value = 0.896551724137931
    $ws.Range("I5").Value = 1
    $ws.Range("J5").Value = 262.5
    $ws.Range("K5").Value = -23.2150423728814
    $ws.Range("L5").Value = -66.775198993745
    $ws.Range("M5").Value = 10.6839047314617
    $ws.Range("N5").Value = -8.84382566585956
    $ws.Range("P5").Value = 'Likely improving'

    # Row 6
    $ws.Range("D6").Value = $true
    $ws.Range("F6").Value = 0.939205768256117
    $ws.Range("G6").Value = 0.0517241379310345
    $ws.Range("H6").Value = 0.982758620689655
    $ws.Range("I6").Value = 2
    $ws.Range("J6").Value = 0.126092739490036
    $ws.Range("K6").Value = -0.0196580376571521
    $ws.Range("L6").Value = -0.0435222417440898
    $ws.Range("M6").Value = -0.0028992497726925
    $ws.Range("N6").Value = -15.5901424115744
    $ws.Range("P6").Value = 'Very likely improving'

    # Row 7
    $ws.Range("F7").Value = 0.971236564272545
    $ws.Range("H7").Value = 0.741379310344828
    $ws.Range("J7").Value = 0.04
    $ws.Range("K7").Value = -0.004969387755102
    $ws.Range("L7").Value = -0.0104357142857143
    $ws.Range("M7").Value = -0.0006941016185985
    $ws.Range("N7").Value = -12.4234693877551
    $ws.Range("P7").Value = 'Extremely likely improving'

    # Row 8
    $ws.Range("F8").Value = 0.97269556738393
    $ws.Range("G8").Value = 0.0172413793103448
    $ws.Range("J8").Value = 0.608
    $ws.Range("K8").Value = -0.0430858093105525
    $ws.Range("L8").Value = -0.0749616218936508
    $ws.Range("M8").Value = -0.0111712340912962
    $ws.Range("N8").Value = -7.08648179449877
    $ws.Range("P8").Value = 'Extremely likely improving'

    # Row 9
    $ws.Range("F9").Value = 0.150355523574156
    $ws.Range("H9").Value = 0.827586206896552
    $ws.Range("J9").Value = 7.865
    $ws.Range("K9").Value = -0.0346520381594968
    $ws.Range("L9").Value = -0.0702891196900986
    $ws.Range("M9").Value = 0.0129446718437319
    $ws.Range("N9").Value = -0.44058535485692
    $ws.Range("P9").Value = 'Unlikely increasing'

    # Row 10
    $ws.Range("D10").Value = $true
    $ws.Range("F10").Value = 0.980765377456021
    $ws.Range("H10").Value = 0.810344827586207
    $ws.Range("J10").Value = 0.88
    $ws.Range("K10").Value = -0.0822124361372456
    $ws.Range("L10").Value = -0.146770709964562
    $ws.Range("M10").Value = -0.0189083199152397
    $ws.Range("N10").Value = -9.34232228832336
    $ws.Range("P10").Value = 'Extremely likely improving'

    # Row 11
    $ws.Range("B11").Value = 'Total Nitrogen'
    $ws.Range("D11").Value = $false
    $ws.Range("E11").Value = 'ok'
    $ws.Range("F11").Value = 0.976730209576497
    $ws.Range("G11").Value = 0
    $ws.Range("H11").Value = 0.913793103448276
    $ws.Range("I11").Value = 0
    $ws.Range("J11").Value = 2.21
    $ws.Range("K11").Value = -0.125474085193157
    $ws.Range("L11").Value = -0.185906084152965
    $ws.Range("M11").Value = -0.0242353145028697
    $ws.Range("N11").Value = -5.67756041597997
    $ws.Range("P11").Value = 'Extremely likely improving'
    $ws.Range("W11").Value = 'g/m3'

    # Row 12
    $ws.Range("B12").Value = 'Total Phosphorus'
    $ws.Range("D12").Value = $true
    $ws.Range("F12").Value = 0.356211124939333
    $ws.Range("H12").Value = 0.948275862068966
    $ws.Range("J12").Value = 0.2525
    $ws.Range("K12").Value = 0.0043322386016901
    $ws.Range("L12").Value = -0.0065502728376941
    $ws.Range("M12").Value = 0.0132223869630501
    $ws.Range("N12").Value = 1.7157380600753
    $ws.Range("W12").Value = 'g/m3'

    # "Suspended Sediment Concentration" (old row 11) is dropped; "Total
    # Nitrogen" and "Total Phosphorus" have shifted up into rows 11/12 above
    # with their refreshed values, so the old trailing row 13 is removed.
    $ws.Rows(13).Delete()
